$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.366.16"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").Value = "2.371.64"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'311.73"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("D6").Value = "'109.00"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("D10").Value = "'40.95"
$ws.Range("E10").Value = "  -0.75%  "
$ws.Range("D11").Value = "'0.0916"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "'8.47"
$ws.Range("E12").Value = "  -1.01%  "
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "'0.977"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("D15").Value = "2.726.86"
$ws.Range("E15").Value = "  +1.29%  "
$ws.Range("D16").Value = "'15.22"
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "2.366.49"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "45.252.94"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "'14.45"
$ws.Range("E19").Value = "  +10.78%  "
$ws.Range("D20").Value = "'7.23"
$ws.Range("E20").Value = "  -4.20%  "
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'73.11"
$ws.Range("E22").Value = "  -1.65%  "
$ws.Range("D23").Value = "'3.50"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").Value = "'259.78"
$ws.Range("E24").Value = "  -3.07%  "
$ws.Range("D25").Value = "'2.31"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'11.10"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'7.22"
$ws.Range("E28").Value = "  -6.26%  "
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("D30").Value = "'0.0969"
$ws.Range("E30").Value = "  +9.10%  "
$ws.Range("D31").Value = "'22.38"
$ws.Range("E31").Value = "  -1.01%  "
$ws.Range("D32").Value = "'37.33"
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("D33").Value = "'168.08"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'3.00"
$ws.Range("E34").Value = "  +5.08%  "
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").Value = "'4.68"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").Value = "'3.98"
$ws.Range("E38").Value = "  +4.47%  "
$ws.Range("D39").Value = "'0.0353"
$ws.Range("E39").Value = "  -3.36%  "
$ws.Range("D40").Value = "'2.90"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("D41").Value = "'1.78"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "'100.04"
$ws.Range("E42").Value = "  -4.63%  "
$ws.Range("D43").Value = "1.903.91"
$ws.Range("E43").Value = "  +14.77%  "
$ws.Range("D44").Value = "'69.66"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "'0.229"
$ws.Range("E45").Value = "  -4.13%  "
$ws.Range("D46").Value = "'12.88"
$ws.Range("E46").Value = "  -5.13%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Value = "'83.48"
$ws.Range("E48").Value = "  +8.91%  "
$ws.Range("D49").Value = "'5.64"
$ws.Range("E49").Value = "  +7.36%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "'9.20"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'110.60"
$ws.Range("E51").Value = "  -3.34%  "
